# Fruta / hortaliza, semanal
# Insert a new weekly record as row 91 (pushing the existing rows 91-150
# down to 92-151) on the single "Coliflor" data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 91..150 down to 92..151, inheriting row 91's existing formatting
# (this is how the source workbook grew from 150 to 151 data rows).
$ws.Rows("91:91").Insert()

# Populate the newly opened row 91 with the new weekly observation.
$ws.Range("A91").Value = 7
$ws.Range("B91").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C91").Value = "Ñuble"
$ws.Range("D91").Value = 44438
$ws.Range("E91").Value = 16
$ws.Range("F91").Value = 100112008
$ws.Range("G91").Value = "Coliflor"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 300
$ws.Range("K91").Value = 700
$ws.Range("L91").Value = 750
$ws.Range("M91").Value = 725
$ws.Range("N91").Value = "$/unidad"
$ws.Range("O91").Value = "Provincia de Diguillín"
$ws.Range("P91").Value = 725
$ws.Range("Q91").Value = 1
$ws.Range("R91").Value = "Hortaliza"
